$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3301.5334
$ws.Range("I76").Value = 3305.4075
$ws.Range("J76").Value = 3266.6667
$ws.Range("K76").Value = 3305.4075
$ws.Range("L76").Value = 3266.6667
$ws.Range("M76").Value = -2990.4075
$ws.Range("N76").Value = -3896.6667
$ws.Range("H79").Value = 3301.5334
$ws.Range("I79").Value = 3305.4075
$ws.Range("J79").Value = 3266.6667
$ws.Range("K79").Value = 3305.4075
$ws.Range("L79").Value = 3266.6667
$ws.Range("M79").Value = -2213.4075
$ws.Range("N79").Value = -5450.6667

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2533.0833
$ws.Range("I63").Value = 2733.2222
$ws.Range("J63").Value = 1932.6666
$ws.Range("K63").Value = 2733.2222
$ws.Range("L63").Value = 1932.6666
$ws.Range("M63").Value = -2047.2222
$ws.Range("N63").Value = -3304.6666
$ws.Range("H66").Value = 2533.0833
$ws.Range("I66").Value = 2733.2222
$ws.Range("J66").Value = 1932.6666
$ws.Range("K66").Value = 13666.111
$ws.Range("L66").Value = 9663.333000000001
$ws.Range("M66").Value = -10234.111
$ws.Range("N66").Value = -16527.333
$ws.Range("H88").Value = 10567696
$ws.Range("J88").Value = 17611858
$ws.Range("L88").Value = 17611858
$ws.Range("N88").Value = -17612670
$ws.Range("H91").Value = 10567696
$ws.Range("J91").Value = 17611858
$ws.Range("L91").Value = 17611858
$ws.Range("N91").Value = -17614666
$ws.Range("H140").Value = 36439.332
$ws.Range("J140").Value = 36439.332
$ws.Range("L140").Value = 36439.332
$ws.Range("N140").Value = -46799.332
$ws.Range("H141").Value = 92659.8
$ws.Range("J141").Value = 92659.8
$ws.Range("L141").Value = 92659.8
$ws.Range("N141").Value = -103019.8

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 23000
$ws.Range("J56").Value = 23000
$ws.Range("L56").Value = 23000
$ws.Range("N56").Value = -24478
$ws.Range("H80").Value = 13860.357
$ws.Range("I80").Value = 27212.428
$ws.Range("J80").Value = 508.2857
$ws.Range("K80").Value = 27212.428
$ws.Range("L80").Value = 508.2857
$ws.Range("M80").Value = -26214.428
$ws.Range("N80").Value = -2504.2857
$ws.Range("H81").Value = 18984.625
$ws.Range("J81").Value = 18984.625
$ws.Range("L81").Value = 18984.625
$ws.Range("N81").Value = -21106.625
$ws.Range("H83").Value = 13860.357
$ws.Range("I83").Value = 27212.428
$ws.Range("J83").Value = 508.2857
$ws.Range("K83").Value = 136062.14
$ws.Range("L83").Value = 2541.4285
$ws.Range("M83").Value = -131070.14
$ws.Range("N83").Value = -12525.4285
$ws.Range("H84").Value = 18984.625
$ws.Range("J84").Value = 18984.625
$ws.Range("L84").Value = 56953.875
$ws.Range("N84").Value = -67561.875
$ws.Range("H105").Value = 886048.5600000001
$ws.Range("I105").Value = 1991009.2
$ws.Range("J105").Value = 2080
$ws.Range("K105").Value = 1991009.2
$ws.Range("L105").Value = 2080
$ws.Range("M105").Value = -1989262.2
$ws.Range("N105").Value = -5574

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7454.643
$ws.Range("I31").Value = 1686.8334
$ws.Range("J31").Value = 11780.5
$ws.Range("K31").Value = 1686.8334
$ws.Range("L31").Value = 11780.5
$ws.Range("M31").Value = -1391.8334
$ws.Range("N31").Value = -12370.5
$ws.Range("H34").Value = 7454.643
$ws.Range("I34").Value = 1686.8334
$ws.Range("J34").Value = 11780.5
$ws.Range("K34").Value = 1686.8334
$ws.Range("L34").Value = 11780.5
$ws.Range("M34").Value = -1484.8334
$ws.Range("N34").Value = -12184.5
$ws.Range("H62").Value = 2490.9167
$ws.Range("I62").Value = 2510.625
$ws.Range("J62").Value = 2451.5
$ws.Range("K62").Value = 2510.625
$ws.Range("L62").Value = 2451.5
$ws.Range("M62").Value = -1886.625
$ws.Range("N62").Value = -3699.5
$ws.Range("H65").Value = 2490.9167
$ws.Range("I65").Value = 2510.625
$ws.Range("J65").Value = 2451.5
$ws.Range("K65").Value = 12553.125
$ws.Range("L65").Value = 12257.5
$ws.Range("M65").Value = -9433.125
$ws.Range("N65").Value = -18497.5
$ws.Range("H68").Value = 17318
$ws.Range("J68").Value = 17318
$ws.Range("L68").Value = 17318
$ws.Range("N68").Value = -18816
$ws.Range("H71").Value = 17318
$ws.Range("J71").Value = 17318
$ws.Range("L71").Value = 51954
$ws.Range("N71").Value = -59442
$ws.Range("H74").Value = 13128.9
$ws.Range("J74").Value = 13128.9
$ws.Range("L74").Value = 13128.9
$ws.Range("N74").Value = -14876.9
$ws.Range("H77").Value = 13128.9
$ws.Range("J77").Value = 13128.9
$ws.Range("L77").Value = 39386.7
$ws.Range("N77").Value = -48122.7
$ws.Range("H86").Value = 43487944
$ws.Range("I86").Value = 83347300
$ws.Range("J86").Value = 5017.8184
$ws.Range("K86").Value = 83347300
$ws.Range("L86").Value = 5017.8184
$ws.Range("M86").Value = -83346177
$ws.Range("N86").Value = -7263.8184
$ws.Range("H89").Value = 43487944
$ws.Range("I89").Value = 83347300
$ws.Range("J89").Value = 5017.8184
$ws.Range("K89").Value = 416736500
$ws.Range("L89").Value = 25089.092
$ws.Range("M89").Value = -416730884
$ws.Range("N89").Value = -36321.092

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 294890.47
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 294890.47
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 884671.4099999999
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -894751.4099999999
$ws.Range("H132").Value = 1588.7037
$ws.Range("I132").Value = 630.3077
$ws.Range("J132").Value = 2478.6428
$ws.Range("K132").Value = 5672.7693
$ws.Range("L132").Value = 22307.7852
$ws.Range("M132").Value = -3142.7693
$ws.Range("N132").Value = -27367.7852

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4200
$ws.Range("I70").Value = 3846.1538
$ws.Range("J70").Value = 6500
$ws.Range("K70").Value = 3846.1538
$ws.Range("L70").Value = 6500
$ws.Range("M70").Value = -3576.1538
$ws.Range("N70").Value = -7040
$ws.Range("H73").Value = 4200
$ws.Range("I73").Value = 3846.1538
$ws.Range("J73").Value = 6500
$ws.Range("K73").Value = 3846.1538
$ws.Range("L73").Value = 6500
$ws.Range("M73").Value = -2910.1538
$ws.Range("N73").Value = -8372
$ws.Range("H80").Value = 8754.0625
$ws.Range("I80").Value = 2937.375
$ws.Range("J80").Value = 14570.75
$ws.Range("K80").Value = 2937.375
$ws.Range("L80").Value = 14570.75
$ws.Range("M80").Value = -1939.375
$ws.Range("N80").Value = -16566.75
$ws.Range("H83").Value = 8754.0625
$ws.Range("I83").Value = 2937.375
$ws.Range("J83").Value = 14570.75
$ws.Range("K83").Value = 14686.875
$ws.Range("L83").Value = 72853.75
$ws.Range("M83").Value = -9694.875
$ws.Range("N83").Value = -82837.75
$ws.Range("H126").Value = 2085.7144
$ws.Range("I126").Value = 2133.3333
$ws.Range("J126").Value = 2077.7778
$ws.Range("K126").Value = 6399.999899999999
$ws.Range("L126").Value = 6233.3334
$ws.Range("M126").Value = -3929.999899999999
$ws.Range("N126").Value = -11173.3334

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1789.4117
$ws.Range("I68").Value = 1272.7273
$ws.Range("J68").Value = 2736.6667
$ws.Range("K68").Value = 1272.7273
$ws.Range("L68").Value = 2736.6667
$ws.Range("M68").Value = -523.7273
$ws.Range("N68").Value = -4234.6667
$ws.Range("H71").Value = 1789.4117
$ws.Range("I71").Value = 1272.7273
$ws.Range("J71").Value = 2736.6667
$ws.Range("K71").Value = 6363.636500000001
$ws.Range("L71").Value = 13683.3335
$ws.Range("M71").Value = -2619.636500000001
$ws.Range("N71").Value = -21171.3335
